$wb = $excel.ActiveWorkbook

# --- TC04 sheet: add Start Date columns ---
$ws4 = $wb.Worksheets.Item("TC04")

$ws4.Range("C1").Value = "startdate"
$ws4.Range("D2").Value = "Start Date"
$ws4.Range("D1").Value = "startcalendartitle"
$ws4.Range("C2").NumberFormat = "mm-dd-yy"
$ws4.Range("C2").Value = (Get-Date -Year 2019 -Month 12 -Day 12 -Hour 0 -Minute 0 -Second 0).Date

# selection changes (active cell) per sheet
$ws4.Range("D4").Select()
